$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> newValue } for column F ("想去人数")
$updates = @{
    "展览" = @{
        2  = 128
        5  = 1495
        9  = 280
        11 = 4793
        12 = 8
        15 = 222
        17 = 164
        21 = 3724
        22 = 631
        23 = 612
        26 = 101
        31 = 568
        32 = 7
        34 = 843
        35 = 2324
        36 = 419
    }
    "全部类型" = @{
        2  = 128
        5  = 1495
        9  = 280
        11 = 4793
        12 = 8
        15 = 222
        17 = 164
        21 = 3724
        22 = 631
        23 = 612
        26 = 101
        31 = 568
        32 = 7
        35 = 843
        36 = 2324
        37 = 419
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
